$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Canada": append row 11 (new data point for date 44105 / Canada)
# ---------------------------------------------------------------------------
$wsCanada = $wb.Worksheets.Item("Canada")

# Copy formatting (incl. date number format / style) down from the last
# existing row so the new row matches the established pattern.
$wsCanada.Range("A10:B10").Copy() | Out-Null
$wsCanada.Range("A11:B11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$wsCanada.Range("A11").Value = 44105
$wsCanada.Range("B11").Value = "Canada"
$wsCanada.Range("C11").Value = 60.9
$wsCanada.Range("D11").Value = 1816.8

# ---------------------------------------------------------------------------
# Sheet "Province": append rows 92-101 (one new date block, all provinces)
# ---------------------------------------------------------------------------
$wsProvince = $wb.Worksheets.Item("Province")

$wsProvince.Range("A82:B91").Copy() | Out-Null
$wsProvince.Range("A92:B101").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$provinceData = @(
    @{Row=92;  Name="Newfoundland & Labrador"; C=15.5;  D=32.7},
    @{Row=93;  Name="Prince Edward Island";     C=19.7;  D=8.5},
    @{Row=94;  Name="Nova Scotia";              C=6.6;   D=43.8},
    @{Row=95;  Name="New Brunswick";            C=27.3;  D=39.6},
    @{Row=96;  Name="Quebec";                   C=53.3;  D=354.2},
    @{Row=97;  Name="Ontario";                  C=81.3;  D=768},
    @{Row=98;  Name="Manitoba";                 C=36.5;  D=49.7},
    @{Row=99;  Name="Saskatchewan";             C=20.1;  D=38.2},
    @{Row=100; Name="Alberta";                  C=56.7;  D=267.1},
    @{Row=101; Name="British Columbia";         C=67.2;  D=215}
)

foreach ($entry in $provinceData) {
    $r = $entry.Row
    $wsProvince.Range("A$r").Value = 44105
    $wsProvince.Range("B$r").Value = $entry.Name
    $wsProvince.Range("C$r").Value = $entry.C
    $wsProvince.Range("D$r").Value = $entry.D
}

# ---------------------------------------------------------------------------
# Mirror the interactive selection state left behind by entering the data.
# ---------------------------------------------------------------------------
$wsCanada.Range("C12").Select() | Out-Null
$wsProvince.Activate() | Out-Null
$wsProvince.Range("C102").Select() | Out-Null
